$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "naoures"
$ws.Range("B6").Value = "naoures"
$ws.Range("A7").Value = "naw"
$ws.Range("B7").Value = "nounou"
